$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 77.13253181199616
$ws.Range("B3").Value = 0.8940851131990224
$ws.Range("B4").Value = 0.06137997724936079
$ws.Range("B5").Value = 0.3849804584980858
